# Weekly update: insert a new Piña price-report row for
# Femacal de La Calera (Coquimbo) and push the existing rows 450-490
# down by one (they all keep their data, just move to row 451-491).
#
# New row 450 data:
#   Fecha=44578, Calidad=Primera, Volumen=160,
#   Precio minimo=16500, Precio maximo=17000, Precio promedio ponderado=16750,
#   Unidad de comercializacion=$/caja 12 unidades, Origen=Ecuador,
#   Precio $/Kg=1396, Kg/unidad=12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 450; Excel shifts rows 450..490 down to 451..491,
# carrying their existing values/formatting along automatically.
$ws.Rows("450:450").Insert()

# Populate the newly-inserted row 450 with the new observation.
$ws.Range("A450").Value = 3
$ws.Range("B450").Value = "Femacal de La Calera"
$ws.Range("C450").Value = "Coquimbo"
$ws.Range("D450").Value = 44578
$ws.Range("E450").Value = 5
$ws.Range("F450").Value = "Fruta"
$ws.Range("G450").Value = 100108
$ws.Range("H450").Value = "Tropicales y subtropicales"
$ws.Range("I450").Value = 100108005
$ws.Range("J450").Value = "Piña"
$ws.Range("K450").Value = "Caramelo"
$ws.Range("L450").Value = "Primera"
$ws.Range("M450").Value = 160
$ws.Range("N450").Value = 16500
$ws.Range("O450").Value = 17000
$ws.Range("P450").Value = 16750
$ws.Range("Q450").Value = "$/caja 12 unidades"
$ws.Range("R450").Value = "Ecuador"
$ws.Range("S450").Value = 1396
$ws.Range("T450").Value = 12
